$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1079.6216
$ws.Range("I28").Value = 935.2273
$ws.Range("J28").Value = 1291.4
$ws.Range("K28").Value = 935.2273
$ws.Range("L28").Value = 1291.4
$ws.Range("M28").Value = -450.2273
$ws.Range("N28").Value = -2261.4
# Row 55
$ws.Range("H55").Value = 116.75
$ws.Range("I55").Value = 80.666664
$ws.Range("J55").Value = 138.4
$ws.Range("K55").Value = 80.666664
$ws.Range("L55").Value = 138.4
$ws.Range("M55").Value = 133.333336
$ws.Range("N55").Value = -566.4
# Row 116
$ws.Range("H116").Value = 1644.4445
$ws.Range("I116").Value = 1326.25
$ws.Range("J116").Value = 1899
$ws.Range("K116").Value = 1326.25
$ws.Range("L116").Value = 1899
$ws.Range("M116").Value = 2115.75
$ws.Range("N116").Value = -8783
# Row 127
$ws.Range("H127").Value = 1264.75
$ws.Range("I127").Value = 517.4
$ws.Range("J127").Value = 1604.4546
$ws.Range("K127").Value = 1552.2
$ws.Range("L127").Value = 4813.3638
$ws.Range("M127").Value = 3407.8
$ws.Range("N127").Value = -14733.3638
# Row 129
$ws.Range("H129").Value = 852.28
$ws.Range("J129").Value = 1005.7143
$ws.Range("L129").Value = 3017.1429
$ws.Range("N129").Value = -13017.1429
# Row 132
$ws.Range("H132").Value = 3064979
$ws.Range("I132").Value = 2520.1538
$ws.Range("K132").Value = 7560.4614
$ws.Range("M132").Value = -5030.4614
# Row 133
$ws.Range("H133").Value = 25434.545
$ws.Range("J133").Value = 25434.545
$ws.Range("L133").Value = 25434.545
$ws.Range("N133").Value = -35554.545
# Row 138
$ws.Range("H138").Value = 2780276.2
$ws.Range("J138").Value = 3790531.5
$ws.Range("L138").Value = 11371594.5
$ws.Range("N138").Value = -11381874.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 6000
$ws.Range("I45").Value = 10000
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -9623
$ws.Range("N45").Value = -2754
# Row 46
$ws.Range("H46").Value = 2417.6
$ws.Range("I46").Value = 2266
$ws.Range("J46").Value = 2645
$ws.Range("K46").Value = 2266
$ws.Range("L46").Value = 2645
$ws.Range("M46").Value = -1947
$ws.Range("N46").Value = -3283

$ws = $wb.Worksheets.Item("BSM")
# Row 109
$ws.Range("H109").Value = 24950
$ws.Range("J109").Value = 24950
$ws.Range("L109").Value = 24950
$ws.Range("N109").Value = -27724

$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 3251.05
$ws.Range("I94").Value = 5262.3
$ws.Range("J94").Value = 1239.8
$ws.Range("K94").Value = 5262.3
$ws.Range("L94").Value = 1239.8
$ws.Range("M94").Value = -4811.3
$ws.Range("N94").Value = -2141.8
# Row 107
$ws.Range("H107").Value = 567.9167
$ws.Range("I107").Value = 540.2857
$ws.Range("J107").Value = 606.6
$ws.Range("K107").Value = 540.2857
$ws.Range("L107").Value = 606.6
$ws.Range("M107").Value = 1379.7143
$ws.Range("N107").Value = -4446.6
# Row 122
$ws.Range("H122").Value = 1932.2106
$ws.Range("I122").Value = 1394.1333
$ws.Range("J122").Value = 3950
$ws.Range("K122").Value = 4182.3999
$ws.Range("L122").Value = 11850
$ws.Range("M122").Value = -1732.3999
$ws.Range("N122").Value = -16750
# Row 132
$ws.Range("H132").Value = 36888.277
$ws.Range("I132").Value = 2240.5217
$ws.Range("J132").Value = 169704.67
$ws.Range("K132").Value = 6721.5651
$ws.Range("L132").Value = 509114.01
$ws.Range("M132").Value = -4191.5651
$ws.Range("N132").Value = -514174.01
# Row 134
$ws.Range("H134").Value = 51637.773
$ws.Range("I134").Value = 3294.1428
$ws.Range("J134").Value = 136239.12
$ws.Range("K134").Value = 9882.428400000001
$ws.Range("L134").Value = 408717.36
$ws.Range("M134").Value = -7347.428400000001
$ws.Range("N134").Value = -413787.36

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 40000216
$ws.Range("I12").Value = 100000140
$ws.Range("J12").Value = 265.86667
$ws.Range("K12").Value = 300000420
$ws.Range("L12").Value = 797.60001
$ws.Range("M12").Value = -300000247
$ws.Range("N12").Value = -1143.60001
# Row 33
$ws.Range("H33").Value = 195.09091
$ws.Range("I33").Value = 110.4
$ws.Range("J33").Value = 265.66666
$ws.Range("K33").Value = 662.4000000000001
$ws.Range("L33").Value = 1593.99996
$ws.Range("M33").Value = -379.4000000000001
$ws.Range("N33").Value = -2159.99996
# Row 92
$ws.Range("H92").Value = 862.5
$ws.Range("I92").Value = 961.1111
$ws.Range("K92").Value = 2883.3333
$ws.Range("M92").Value = -1635.3333
# Row 97
$ws.Range("H97").Value = 2565.5557
$ws.Range("I97").Value = 3681
$ws.Range("J97").Value = 334.66666
$ws.Range("K97").Value = 11043
$ws.Range("L97").Value = 1003.99998
$ws.Range("M97").Value = -10547
$ws.Range("N97").Value = -1995.99998
# Row 107
$ws.Range("H107").Value = 420.54166
$ws.Range("I107").Value = 689.5625
$ws.Range("J107").Value = 286.03125
$ws.Range("K107").Value = 2068.6875
$ws.Range("L107").Value = 858.09375
$ws.Range("M107").Value = -148.6875
$ws.Range("N107").Value = -4698.09375
# Row 113
$ws.Range("H113").Value = 641.5599999999999
$ws.Range("J113").Value = 663.7
$ws.Range("L113").Value = 1991.1
$ws.Range("N113").Value = -6331.1
# Row 131
$ws.Range("H131").Value = 1007.92957
$ws.Range("I131").Value = 687
$ws.Range("J131").Value = 1043.0312
$ws.Range("K131").Value = 2061
$ws.Range("L131").Value = 3129.0936
$ws.Range("M131").Value = 2979
$ws.Range("N131").Value = -13209.0936

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1183.4286
$ws.Range("I102").Value = 1050
$ws.Range("J102").Value = 1361.3334
$ws.Range("K102").Value = 1050
$ws.Range("L102").Value = 1361.3334
$ws.Range("M102").Value = 572
$ws.Range("N102").Value = -4605.3334
# Row 122
$ws.Range("H122").Value = 2411.524
$ws.Range("I122").Value = 1767.6428
$ws.Range("J122").Value = 3699.2856
$ws.Range("K122").Value = 5302.928400000001
$ws.Range("L122").Value = 11097.8568
$ws.Range("M122").Value = -2852.928400000001
$ws.Range("N122").Value = -15997.8568
# Row 126
$ws.Range("H126").Value = 1863.4517
$ws.Range("I126").Value = 1106.1666
$ws.Range("J126").Value = 2912
$ws.Range("K126").Value = 3318.4998
$ws.Range("L126").Value = 8736
$ws.Range("M126").Value = -848.4998000000001
$ws.Range("N126").Value = -13676

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1070.2778
$ws.Range("I22").Value = 677.1429000000001
$ws.Range("J22").Value = 1320.4546
$ws.Range("K22").Value = 677.1429000000001
$ws.Range("L22").Value = 1320.4546
$ws.Range("M22").Value = -382.1429000000001
$ws.Range("N22").Value = -1910.4546
# Row 27
$ws.Range("H27").Value = 1070.2778
$ws.Range("I27").Value = 677.1429000000001
$ws.Range("J27").Value = 1320.4546
$ws.Range("K27").Value = 677.1429000000001
$ws.Range("L27").Value = 1320.4546
$ws.Range("M27").Value = -570.1429000000001
$ws.Range("N27").Value = -1534.4546
# Row 46
$ws.Range("H46").Value = 748
$ws.Range("I46").Value = 745
$ws.Range("J46").Value = 751
$ws.Range("K46").Value = 745
$ws.Range("L46").Value = 751
$ws.Range("M46").Value = -557
$ws.Range("N46").Value = -1127
# Row 61
$ws.Range("H61").Value = 2754.5334
$ws.Range("I61").Value = 2408.3076
$ws.Range("K61").Value = 2408.3076
$ws.Range("M61").Value = -2206.3076
# Row 93
$ws.Range("H93").Value = 1246.1578
$ws.Range("I93").Value = 1341.0834
$ws.Range("J93").Value = 1083.4286
$ws.Range("K93").Value = 1341.0834
$ws.Range("L93").Value = 1083.4286
$ws.Range("M93").Value = -93.08339999999998
$ws.Range("N93").Value = -3579.4286
# Row 113
$ws.Range("H113").Value = 2754.5334
$ws.Range("I113").Value = 2408.3076
$ws.Range("K113").Value = 2408.3076
$ws.Range("M113").Value = -238.3076000000001
# Row 132
$ws.Range("H132").Value = 153629.8
$ws.Range("I132").Value = 102760
$ws.Range("J132").Value = 204499.6
$ws.Range("K132").Value = 308280
$ws.Range("L132").Value = 613498.8
$ws.Range("M132").Value = -305750
$ws.Range("N132").Value = -618558.8
# Row 136
$ws.Range("H136").Value = 286380.72
$ws.Range("I136").Value = 167444.17
$ws.Range("K136").Value = 502332.51
$ws.Range("M136").Value = -499782.51

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 264.80768
$ws.Range("I107").Value = 265.07693
$ws.Range("J107").Value = 264.53845
$ws.Range("K107").Value = 795.2307900000001
$ws.Range("L107").Value = 793.61535
$ws.Range("M107").Value = 1124.76921
$ws.Range("N107").Value = -4633.61535
# Row 122
$ws.Range("H122").Value = 4133.3335
$ws.Range("I122").Value = 1600
$ws.Range("J122").Value = 4640
$ws.Range("K122").Value = 4800
$ws.Range("L122").Value = 13920
$ws.Range("M122").Value = -2350
$ws.Range("N122").Value = -18820
# Row 129
$ws.Range("H129").Value = 28910
$ws.Range("J129").Value = 28910
$ws.Range("L129").Value = 28910
$ws.Range("N129").Value = -38910
